$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.930.56"
$ws.Range("E2").Value = "  +4.43%  "
$ws.Range("D3").Value = "2.644.66"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.175"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.92%  "
$ws.Range("D10").Value = "2.642.77"
$ws.Range("E10").Value = "  +3.05%  "
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000190"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.28%  "
$ws.Range("D15").Value = "3.124.95"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "72.667.04"
$ws.Range("E16").Value = "  +4.19%  "
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "2.635.57"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "385.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("E23").Value = "  +15.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.17%  "
$ws.Range("D28").Value = "2.778.08"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "0.0₃0963"
$ws.Range("E30").Value = "  +4.09%  "
$ws.Range("E31").Value = "  +4.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "521.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  +4.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.113"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  +5.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.548"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.41%  "
$ws.Range("E50").Value = "  +4.55%  "
$ws.Range("E51").Value = "  +3.03%  "
